# Unit10_Assertion.pptx edits
#
# 1) Slide 3 (sldId 552): Content Placeholder 5 - expand the "Can be written
#    anywhere..." bullet with a parenthetical about assert() / assert.h, and
#    move TextBox 15 up a bit (y offset change only).
# 2) Slide 4 (sldId 553): Content Placeholder 5 - "an assertion" -> "assertions".
# 3) Slide 5 (sldId 554): Content Placeholder 5 - same wording fix, split
#    across three runs as the author re-typed the middle of the sentence.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 (index 3 of the deck) - sldId 552
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Content Placeholder 5 -> 3rd bullet paragraph ("Can be written anywhere...")
$shape3Content = $s3.Shapes.Item(4)
$tr3 = $shape3Content.TextFrame.TextRange
$para = $tr3.Paragraphs(3, 1)
$run1 = $para.Runs(1)

$partA = "Can be written anywhere in the code and checked automatically (e.g., using assert() from "
$partB = "assert.h"
$partC = ")"

# Extend the first run's text in place (keeps its original run/formatting).
$run1.Text = $partA + $partB + $partC

# Carve the tail into two more runs without touching their formatting by
# reassigning each character range's Text to itself - this creates a clean
# run boundary while leaving the inherited rPr untouched.
$startB = $partA.Length + 1
$startC = $partA.Length + $partB.Length + 1

$runB = $para.Characters($startB, $partB.Length)
$runB.Text = $partB

$runC = $para.Characters($startC, $partC.Length)
$runC.Text = $partC

# TextBox 15 - move up (only the vertical offset changes; 3970726 EMU).
$shape3Box = $s3.Shapes.Item(5)
$shape3Box.Top = (3970726 / 12700) + 0.00002

# ---------------------------------------------------------------------------
# Slide 4 (index 4 of the deck) - sldId 553
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shape4Content = $s4.Shapes.Item(4)
$tr4 = $shape4Content.TextFrame.TextRange
$para4 = $tr4.Paragraphs(1, 1)
$run4 = $para4.Runs(1)
$run4.Text = "Important: Write assertions based on "

# ---------------------------------------------------------------------------
# Slide 5 (index 5 of the deck) - sldId 554
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shape5Content = $s5.Shapes.Item(1)
$tr5 = $shape5Content.TextFrame.TextRange
$para5 = $tr5.Paragraphs(1, 1)
$run5 = $para5.Runs(1)

$partA5 = "Important"
$partB5 = ": Write assertions "
$partC5 = "based on your "

$run5.Text = $partA5 + $partB5 + $partC5

$startB5 = $partA5.Length + 1
$startC5 = $partA5.Length + $partB5.Length + 1

$runB5 = $para5.Characters($startB5, $partB5.Length)
$runB5.Text = $partB5

$runC5 = $para5.Characters($startC5, $partC5.Length)
$runC5.Text = $partC5
